$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 162.26666
$ws.Range("I11").Value = 162.26666
$ws.Range("K11").Value = 162.26666
$ws.Range("M11").Value = -22.26666

$ws.Range("H33").Value = 442.09525
$ws.Range("J33").Value = 475.33334
$ws.Range("L33").Value = 475.33334
$ws.Range("N33").Value = -933.33334

$ws.Range("H106").Value = 3303.04
$ws.Range("I106").Value = 1723.0625
$ws.Range("J106").Value = 6111.8887
$ws.Range("K106").Value = 1723.0625
$ws.Range("L106").Value = 6111.8887
$ws.Range("M106").Value = -1092.0625
$ws.Range("N106").Value = -7373.8887

$ws.Range("H111").Value = 979.4286
$ws.Range("I111").Value = 979.4286
$ws.Range("K111").Value = 2938.2858
$ws.Range("M111").Value = 128.7142000000003

$ws.Range("H135").Value = 815.9231
$ws.Range("I135").Value = 815.9231
$ws.Range("K135").Value = 7343.3079
$ws.Range("M135").Value = -4808.3079

$ws.Range("H137").Value = 3073.476
$ws.Range("J137").Value = 3099.6182
$ws.Range("L137").Value = 9298.854599999999
$ws.Range("N137").Value = -14398.8546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4215.216
$ws.Range("I32").Value = 3020.4167
$ws.Range("K32").Value = 3020.4167
$ws.Range("M32").Value = -2733.4167

$ws.Range("H74").Value = 12347954
$ws.Range("I74").Value = 13334898
$ws.Range("K74").Value = 13334898
$ws.Range("M74").Value = -13334024

$ws.Range("H77").Value = 12347954
$ws.Range("I77").Value = 13334898
$ws.Range("K77").Value = 66674490
$ws.Range("M77").Value = -66670122

$ws.Range("H102").Value = 2637.7144
$ws.Range("I102").Value = 1243.4
$ws.Range("J102").Value = 6123.5
$ws.Range("K102").Value = 1243.4
$ws.Range("L102").Value = 6123.5
$ws.Range("M102").Value = 378.5999999999999
$ws.Range("N102").Value = -9367.5

$ws.Range("H132").Value = 2561.8708
$ws.Range("I132").Value = 1770.16
$ws.Range("J132").Value = 5860.6665
$ws.Range("K132").Value = 5310.48
$ws.Range("L132").Value = 17581.9995
$ws.Range("M132").Value = -2780.48
$ws.Range("N132").Value = -22641.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 699.5
$ws.Range("J64").Value = 1000
$ws.Range("L64").Value = 1000
$ws.Range("N64").Value = -1450

$ws.Range("H67").Value = 699.5
$ws.Range("J67").Value = 1000
$ws.Range("L67").Value = 1000
$ws.Range("N67").Value = -2560

$ws.Range("H86").Value = 6047.6665
$ws.Range("J86").Value = 7612.6665
$ws.Range("L86").Value = 7612.6665
$ws.Range("N86").Value = -9858.666499999999

$ws.Range("H89").Value = 6047.6665
$ws.Range("J89").Value = 7612.6665
$ws.Range("L89").Value = 38063.3325
$ws.Range("N89").Value = -49295.3325

$ws.Range("H107").Value = 1980.2307
$ws.Range("I107").Value = 1869.091
$ws.Range("J107").Value = 2591.5
$ws.Range("K107").Value = 1869.091
$ws.Range("L107").Value = 2591.5
$ws.Range("M107").Value = 50.90900000000011
$ws.Range("N107").Value = -6431.5

$ws.Range("H130").Value = 60000
$ws.Range("J130").Value = 60000
$ws.Range("L130").Value = 60000
$ws.Range("N130").Value = -70040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22131.094
$ws.Range("I31").Value = 1874.8
$ws.Range("J31").Value = 84458.16
$ws.Range("K31").Value = 1874.8
$ws.Range("L31").Value = 84458.16
$ws.Range("M31").Value = -1579.8
$ws.Range("N31").Value = -85048.16

$ws.Range("H34").Value = 22131.094
$ws.Range("I34").Value = 1874.8
$ws.Range("J34").Value = 84458.16
$ws.Range("K34").Value = 1874.8
$ws.Range("L34").Value = 84458.16
$ws.Range("M34").Value = -1672.8
$ws.Range("N34").Value = -84862.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2396
$ws.Range("I22").Value = 744.5
$ws.Range("K22").Value = 2233.5
$ws.Range("M22").Value = -2064.5

$ws.Range("H27").Value = 2396
$ws.Range("I27").Value = 744.5
$ws.Range("K27").Value = 2233.5
$ws.Range("M27").Value = -2131.5

$ws.Range("H133").Value = 5355.25
$ws.Range("I133").Value = 1436.25
$ws.Range("J133").Value = 9274.25
$ws.Range("K133").Value = 4308.75
$ws.Range("L133").Value = 27822.75
$ws.Range("M133").Value = 751.25
$ws.Range("N133").Value = -37942.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 338019.66
$ws.Range("I80").Value = 1251571
$ws.Range("J80").Value = 5819.1816
$ws.Range("K80").Value = 1251571
$ws.Range("L80").Value = 5819.1816
$ws.Range("M80").Value = -1250573
$ws.Range("N80").Value = -7815.1816

$ws.Range("H83").Value = 338019.66
$ws.Range("I83").Value = 1251571
$ws.Range("J83").Value = 5819.1816
$ws.Range("K83").Value = 6257855
$ws.Range("L83").Value = 29095.908
$ws.Range("M83").Value = -6252863
$ws.Range("N83").Value = -39079.908

$ws.Range("H97").Value = 2738
$ws.Range("I97").Value = 2291
$ws.Range("J97").Value = 3706.5
$ws.Range("K97").Value = 2291
$ws.Range("L97").Value = 3706.5
$ws.Range("M97").Value = -1795
$ws.Range("N97").Value = -4698.5

$ws.Range("H126").Value = 4437.875
$ws.Range("J126").Value = 5851.4
$ws.Range("L126").Value = 17554.2
$ws.Range("N126").Value = -22494.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9553.817999999999
$ws.Range("I7").Value = 5588
$ws.Range("K7").Value = 5588
$ws.Range("M7").Value = -5476

$ws.Range("H40").Value = 10200.904
$ws.Range("I40").Value = 9747.362999999999
$ws.Range("K40").Value = 9747.362999999999
$ws.Range("M40").Value = -9611.362999999999

$ws.Range("H122").Value = 128915.09
$ws.Range("I122").Value = 185155.9
$ws.Range("J122").Value = 5185.3
$ws.Range("K122").Value = 555467.7
$ws.Range("L122").Value = 15555.9
$ws.Range("M122").Value = -553017.7
$ws.Range("N122").Value = -20455.9

$ws.Range("H126").Value = 9553.817999999999
$ws.Range("I126").Value = 5588
$ws.Range("K126").Value = 16764
$ws.Range("M126").Value = -14294

$ws.Range("H132").Value = 4536.885
$ws.Range("I132").Value = 4059.4443
$ws.Range("K132").Value = 12178.3329
$ws.Range("M132").Value = -9648.332900000001

$ws.Range("H136").Value = 2847.2769
$ws.Range("I136").Value = 2161.48
$ws.Range("J136").Value = 5133.2666
$ws.Range("K136").Value = 6484.440000000001
$ws.Range("L136").Value = 15399.7998
$ws.Range("M136").Value = -3934.440000000001
$ws.Range("N136").Value = -20499.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3557.25
$ws.Range("I122").Value = 1897.75
$ws.Range("J122").Value = 6876.25
$ws.Range("K122").Value = 5693.25
$ws.Range("L122").Value = 20628.75
$ws.Range("M122").Value = -3243.25
$ws.Range("N122").Value = -25528.75

$ws.Range("H132").Value = 1993.4348
$ws.Range("J132").Value = 3874.75
$ws.Range("L132").Value = 11624.25
$ws.Range("N132").Value = -16684.25

$ws.Range("H135").Value = 69799.336
$ws.Range("J135").Value = 69799.336
$ws.Range("L135").Value = 69799.336
$ws.Range("N135").Value = -79939.336
